# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.985.13"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "3.149.42"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'574.90"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("D6").Value = "'149.32"
$ws.Range("E6").Value = "  +4.58%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.149.91"
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("D11").Value = "'6.13"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "'0.497"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("E13").Value = "  +12.64%  "
$ws.Range("D14").Value = "'36.99"
$ws.Range("D15").Value = "3.666.79"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").Value = "65.061.49"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "3.156.84"
$ws.Range("E17").Value = "  +3.02%  "
$ws.Range("D18").Value = "'7.09"
$ws.Range("E18").Value = "  +4.45%  "
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "'505.51"
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("D21").Value = "'14.75"
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("D23").Value = "'15.28"
$ws.Range("E23").Value = "  +3.92%  "
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "'84.11"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +3.46%  "
$ws.Range("D28").Value = "'8.83"
$ws.Range("E28").Value = "  +7.74%  "
$ws.Range("E29").Value = "  +5.15%  "
$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").Value = "'2.79"
$ws.Range("E30").Value = "  +7.96%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'27.56"
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'1.18"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").Value = "'6.18"
$ws.Range("E34").Value = "  +8.42%  "
$ws.Range("E35").Value = "  +4.11%  "
$ws.Range("D36").Value = "'54.93"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +9.76%  "
$ws.Range("D38").Value = "'463.35"
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").Value = "'2.98"
$ws.Range("E40").Value = "  +7.43%  "
$ws.Range("D41").Value = "'8.65"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("D42").Value = "3.050.27"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("E44").Value = "  +7.62%  "
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").Value = "'28.53"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").Value = "0.0₃0580"
$ws.Range("E47").Value = "  +11.77%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "'2.24"
$ws.Range("E50").Value = "  +4.89%  "
$ws.Range("D51").Value = "'119.44"
$ws.Range("E51").Value = "  +1.23%  "
